$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number must be forced to Text format
# first, otherwise Excel auto-converts the numeric-looking string to a real number
# (the source data models these as text, e.g. "2.40" must keep its trailing zero).

$ws.Range("D2").Value = "36.436.60"
$ws.Range("E2").Value = "  -2.65%  "
$ws.Range("D3").Value = "1.990.09"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.49"
$ws.Range("E5").Value = "  -8.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.597"
$ws.Range("E6").Value = "  -3.42%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.49"
$ws.Range("E8").Value = "  -3.05%  "
$ws.Range("E9").Value = "  -4.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.18"
$ws.Range("E10").Value = "  +2.57%  "
$ws.Range("E11").Value = "  -3.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0984"
$ws.Range("E12").Value = "  -3.38%  "
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").Value = "2.282.01"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.21"
$ws.Range("E15").Value = "  -3.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.755"
$ws.Range("E16").Value = "  -6.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.07"
$ws.Range("E17").Value = "  -3.44%  "
$ws.Range("D18").Value = "1.961.97"
$ws.Range("E18").Value = "  -3.29%  "
$ws.Range("D19").Value = "36.446.47"
$ws.Range("E19").Value = "  -2.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.85"
$ws.Range("D21").Value = "0.0₃0803"
$ws.Range("E21").Value = "  -4.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.25"
$ws.Range("E22").Value = "  +1.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "221.19"
$ws.Range("E23").Value = "  -3.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.40"
$ws.Range("E26").Value = "  -10.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.98"
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.67"
$ws.Range("E28").Value = "  -3.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.128"
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.90"
$ws.Range("E30").Value = "  -4.19%  "
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.117"
$ws.Range("E32").Value = "  -2.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.38"
$ws.Range("E33").Value = "  -5.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0606"
$ws.Range("E34").Value = "  -6.49%  "
$ws.Range("E35").Value = "  -6.79%  "
$ws.Range("E36").Value = "  -2.68%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.77"
$ws.Range("E38").Value = "  -3.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.29"
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.48"
$ws.Range("E40").Value = "  +5.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.99"
$ws.Range("E41").Value = "  -1.65%  "
$ws.Range("D42").Value = "1.454.41"
$ws.Range("E42").Value = "  +4.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0913"
$ws.Range("E43").Value = "  -2.79%  "
$ws.Range("E44").Value = "  -5.56%  "
$ws.Range("E45").Value = "  -9.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.18"
$ws.Range("E46").Value = "  -1.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "14.97"
$ws.Range("E47").Value = "  -4.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.996"
$ws.Range("E48").Value = "  -2.75%  "
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.82"
$ws.Range("E50").Value = "  -3.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.63"
$ws.Range("E51").Value = "  +17.74%  "
